$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "templateCode" column (column A). This shifts lotNumber,
# catalogNumber, activity, concentration, purity, formulation one
# column to the left (B->A ... G->F).
$ws.Range("A1").EntireColumn.Delete()

# Matches the saved selection state recorded in the diff.
[void]$ws.Range("H7").Select()
